# Add a new "BAT/Volt" row to the escData voltage plot (Plot No 3),
# pushing the existing ESC1..ESC4 current rows (Plot No 4) down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 18 - Excel shifts rows 18:21 down to 19:22 and
# copies formatting (incl. the centered "S" style in column D) from the
# row above, same as a manual Insert would.
$ws.Rows("18:18").Insert()

# Populate the new voltage row for the battery channel.
$ws.Range("A18").Value2 = 3
$ws.Range("B18").Value2 = 2
$ws.Range("C18").Value2 = 2
$ws.Range("D18").Value2 = "S"
$ws.Range("E18").Value2 = "Time [ s ]"
$ws.Range("F18").Value2 = "Voltage~[~V~]"
$ws.Range("G18").Value2 = "Vertical"
$ws.Range("K18").Value2 = "BAT/Volt"
$ws.Range("O18").Value2 = 1
$ws.Range("R18").Value2 = "BATT"

# Restore the selection to the newly added row/column, matching the
# author's final cursor position in the source file.
$ws.Range("R18").Select()
